$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.235.57'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.656.10'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  -0.57%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '219.29'
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5241'
$ws.Range("E6").Value = '  -1.35%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.005'
$ws.Range("E7").Value = '  -0.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2665'
$ws.Range("E8").Value = '  +0.57%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06357'
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.68'
$ws.Range("E10").Value = '  -0.69%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07725'
$ws.Range("E11").Value = '  -1.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.593'
$ws.Range("E12").Value = '  +1.85%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.651.15'
$ws.Range("E13").Value = '  -1.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.884.38'
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5629'
$ws.Range("E15").Value = '  +0.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8248'
$ws.Range("E16").Value = '  +1.20%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.41'
$ws.Range("E17").Value = '  -0.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.235.09'
$ws.Range("E18").Value = '  -0.38%  '
$ws.Range("E19").Value = '  -0.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.695'
$ws.Range("E20").Value = '  -0.20%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.49'
$ws.Range("E21").Value = '  -2.16%  '
$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.42'
$ws.Range("E22").Value = '  +1.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.995'
$ws.Range("E23").Value = '  -0.73%  '
$ws.Range("E24").Value = '  -0.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.75'
$ws.Range("E25").Value = '  -1.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1204'
$ws.Range("E26").Value = '  -1.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.276'
$ws.Range("E27").Value = '  +0.60%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.93'
$ws.Range("E28").Value = '  -1.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.515'
$ws.Range("E29").Value = '  +0.23%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05643'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.278'
$ws.Range("E31").Value = '  -0.44%  '
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.360'
$ws.Range("E33").Value = '  +1.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.583'
$ws.Range("E34").Value = '  -1.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.803'
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9526'
$ws.Range("E36").Value = '  -1.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.415'
$ws.Range("E37").Value = '  -0.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5755'
$ws.Range("E38").Value = '  -0.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01599'
$ws.Range("E39").Value = '  -0.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.002'
$ws.Range("E40").Value = '  +1.09%  '
$ws.Range("E41").Value = '  -1.15%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.004'
$ws.Range("E42").Value = '  -0.59%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8431'
$ws.Range("E43").Value = '  -2.30%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.91'
$ws.Range("E44").Value = '  -0.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.012.44'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.795.37'
$ws.Range("E46").Value = '  -0.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '58.40'
$ws.Range("E47").Value = '  +0.24%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₈107'
$ws.Range("E48").Value = '  -1.34%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05354'
$ws.Range("E49").Value = '  +3.91%  '
$ws.Range("B50").Value = 'Frax'
$ws.Range("C50").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.006'
$ws.Range("E50").Value = '  -0.88%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.046'
$ws.Range("E51").Value = '  +0.70%  '
